$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the sheet ("My Series" -> "Data")
$ws.Name = "Data"

# 2) Update the comment text on A1 to the new CEIC add-in metadata blob
$comment = $ws.Range("A1").Comment
$comment.Text("BRoAAB+LCAAAAAAAAAOlWVtvG8cV/isLPrVAyb3QsiVhvAFvcoiSokBSlZWXYrg7Eqda7jK7s5L4lgItUqQpiqJwivSKPqUoUNdoEyC1e/kvgSW7T/0LPXPZ2QupuHQNwdr5zmXOnDlzzpkReud6ERiXJE5oFD6s2Q2rZpDQi3wanj+speysbt+vveOi3rVHgiMc4wVhwGyAVJjsXyf0YW3O2HLfNK+urhpXzUYUn5uOZdnm4+Fg4s3JAtdpmDAceqSmpfw3S9Vc1PEXQ8KwjxmWkg9r/Um/0SHU6wI2xCE+J3GjnSY0JEnSCxlllCRcMiaYkU53+B25MNdp3G/YyFzDc852SgNf8pU4Ja74YFoypQviOpa9W7d26017ajv7zs6+4zQsa++9TFAzogFO2ITEl9QTwIThxVKIW7tN23acHcdB5kYm0JU7wEWjwB+TS5oQv0OCINnKI6bawJbHYNXbOdNCZkFWKXp7Ex7FeDmfUhaQ7cwYD9vGIlS25EpcdBDFxAP/vZVJh+RqFCu3TpcDoE7nNGarLl5tres4IfFoyZ20naiLulHIWgGJ2fES9pr4EApAcFmcEmTeQcyFujTx4JuGKfHdMxwkRaESEZ1E8UWyxB45hHNsch1XYRBhHwKO0YRRL590jYCO4mgJGmHydhT4B6BVMW8gaM39EFzMp21H0UVu3SYiErsq9hf2dIFZxr6Go8k8uhqFwWqSzhIvpjPid9sZ90Ya4gdSSXfShEULsCKHkMQKyAr+wQGswqhLPLrAwVEATkzcJmgpAaiVsuiMsk4UpIswyWyqoOgEVjQl13qFeoxGsLkhd3oU9sOKlo20ssQ4usq2cB0XTijArcTLtnudUGXuApZt3zpF7Ahf5QENoD4U96KAlqNiMieEbQwJSUE8FR7wiuO2V3xOZOYIgsiE6AbUtff2HtQtG36mlrUvfmBmTUa90BcfUF/u122H5+ycLyOiw3QxmsEJvhRrcm2gVSAEqwjaAQ4vAD2hbH7YyqzfQEFyzXfyr9MQnNVlgFcC1n4pYqgfekHqE5kC+uGZCEpum+K+k47WoAEcaxfhcDVdLXkqMO/g6JIznAZQnRgkkXMdi1UYtZKLKk8RQsdxkO2xy2t/AsXf8xcND9IjL3ANL1pwwISaezJBZpGf53iP9MLzAQ7PU8ii2o9VXEcYzwDTGIcJX45OmpVg28yEsn2Rydxd31xkVljQlCyWUYyDIfiEHqShqJWqGIBzh5jN1QjObUC8zL9mLqqlykZlNr+JTZxAuQKe3NXBroCCice77DAKpzQHEV/jMPKhyOOAzuJSgG2kwWblpS8LNb66Lctg5n9oMyG3fJuseO+RDxTO49W1M4IYoQ6Y5E7G93adHasJ/ZQYI7HkMcGB0YMGlBGjH16ShC1AbN8Yk4T68EVxsG+8S2aEwoEXPlJHbmvpohw6iMn7KTTSK2FKC+wtI2UGyKbnNMTBOqOm5ALuKcFxsCowyqUOIg/4bn/0r5tfP3v5/LPbj568/uIH//n7L1/+46c3T38IH7d/+evNx7+Qy5TMaIpnAREGTdu7u1bzHgSahpBKDVDf/dRjAjs9FfVej5FqW8Wg0+t3Hg3aIpdoMBPnrUTKD8MAr6I0H07kIsREYkvNLBIkizvNcpMal6hdwou8PIW8Q70kZe4i/S5B6YtXLz579eJPd0orh5XKzg4vJ28uO7a9xqfLDu/8ZWOX3Svu1a2duuMUmCs8aAzXI+jjtJ/6vtu0rT3LaVq2zuO+DuRNTFWS0jTF52ZFTkKdKA1ZvNIhUBxnRBH4UzgimiyPQmGgQvSLH7/+85MSl/KuQspawLgojWU0mdlAqD4cT43J6Hjc6RnT3oTHSU4r8EnlX8OsZtfnqRRUYZji4FtwKeZXYqMG/V7NiM4Mgr25sYKTWDiHpWDbhMqJ3lJl1cpHcZQu5Y4UBHJ0A6fOJhslNuQaQRP+XEs6OWkDu7T15m+fbxJQC1FhdhxSpm93RQyVKBIq0NWp/fSfL7/88OXz57fPfnbz5fdLGtQ8+tIDcQ6nqTjUYQ8pT9WbCoJOJsKZF9Z3C/VFgbxpPIpoyBLXvi/6RTVCIGpzbeI36i+g5AnFwl+AVxD0Lk5610wdbPcQmWUA7FxiqLZR3llrQObw3K///s1vb3/1+e0nz15/+Mebj/5w8/Enr1787vXT38tTd/vk2e1PnqosXy0Ewhber8sG0BC3P8/gp9Hgtdv46oOfG2HEDOg5jFRkpK8++LSgjBsqupNcM7Rz2pCyCWusRWEuZxRM0TaU5LSIbAA6vIQ1NYcqYtGSevkk79W5Kn7uBOEb/Wk9TYgRQTv1TVhJmTkX/l/llIgsqUcPLMd2FFVaw5cww0nB9Y+CaAZNRkYQ16sKS0nq6wVyXjHfo8Go3RrkLNKIUezDbczid0T+gbKekpeUfpKN9KUiR4AKjZ+XBvxGvMa2TtKaC2nMVFfLs5bP059rWeLuVq6PJQ64ncexbIhC9VI5SZfQDWcPEHfTxaNMoQE+lL1qsSXOx/1umQ7jAhUKYZnMAUEXqUmRZJrqJ/wWK9vZQ+6afAi00kMOuEM9RspO6xL6ytjkeacXx1G8MfnklIxtCJ00ZBQz97jmEXsqu24/36sMyBKe/pC3PrXCqEsCwrZ7qTNz6WF0+daysPfbivaTUeArZ2539dBuyRUUnyt5oPy/r5Uy2FpxDI0Vf97Y+nkxu7SO4a67pTVyKUKQXwFhdvX+d0DjhD3mmUB9SeRUI6eyQ30MjansQx9LQHDID7VIpd0smZkdXSYfjaNgQBd0y2uhlZ3vshLw5XIpW7j+dpHCS8shuYYGs6ABkuLse1A2+JVnO20yYCGXann+TpPQ8znb1rAHM0x8MrPq3ow49Xu+tVvfI6RZt234H3uOY1k7/JFHKYfMQcnVlpOY2Yblf8xx/ws+naRuBRoAAA==")

# 3) Update B1:K1 values
$ws.Range("B1").Value = 15733.82555555556
$ws.Range("C1").Value = 139798646.1781028
$ws.Range("D1").Value = 11823.64775262282
$ws.Range("E1").Value = 0.951248095277063
$ws.Range("F1").Value = -0.1855642624124751
$ws.Range("G1").Value = 0.7514795248538862
$ws.Range("H1").Value = 3271.34
$ws.Range("I1").Value = 37962.97
$ws.Range("J1").Value = 10540.67
$ws.Range("K1").Value = 9
$ws.Range("L1").Value = 3271.34

# 4) Widen the number format applied to L1 from "0.000" to "###0.000" first
#    (while it is still the sole cell using that style), then copy the
#    resulting format onto the new M1:U1 cells so the whole L1:U1 block
#    shares one consistent style/number format.
$ws.Range("L1").NumberFormat = "###0.000"
$ws.Range("L1").Copy()
$ws.Range("M1:U1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 5) New formula + values for M1:U1
$ws.Range("M1").Formula = "=NA()"
$ws.Range("N1").Value = 6384.48
$ws.Range("O1").Value = 7111.74
$ws.Range("P1").Value = 8354.02
$ws.Range("Q1").Value = 10540.67
$ws.Range("R1").Value = 16421.36
$ws.Range("S1").Value = 22352.63
$ws.Range("T1").Value = 29205.22
$ws.Range("U1").Value = 37962.97
